# Atualização automática de preços de eletricidade
# Updates row 2 of the Spot_PT sheet with the new daily spot prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spot_PT")

# Day (date serial number)
$ws.Range("A2").Value = 45996

# Hourly prices 0h-1h .. 23h-24h (columns B..Y)
$ws.Range("B2").Value = 73.89
$ws.Range("C2").Value = 62.56
$ws.Range("D2").Value = 56.61
$ws.Range("E2").Value = 55.32
$ws.Range("F2").Value = 54.75
$ws.Range("G2").Value = 59.44
$ws.Range("H2").Value = 71.08
$ws.Range("I2").Value = 79.89
$ws.Range("J2").Value = 88.01000000000001
$ws.Range("K2").Value = 84.84999999999999
$ws.Range("L2").Value = 82.03
$ws.Range("M2").Value = 74.78
$ws.Range("N2").Value = 65.23999999999999
$ws.Range("O2").Value = 55.42
$ws.Range("P2").Value = 54.84
$ws.Range("Q2").Value = 67.88
$ws.Range("R2").Value = 78.29000000000001
$ws.Range("S2").Value = 84.31999999999999
$ws.Range("T2").Value = 87.16
$ws.Range("U2").Value = 83
$ws.Range("V2").Value = 79.19
$ws.Range("W2").Value = 71.81999999999999
$ws.Range("X2").Value = 67.94
$ws.Range("Y2").Value = 56.62

# Daily average price
$ws.Range("Z2").Value = 70.62

# Slot_4h_max / Slot_4h_price
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 83.19

# Slot_2h_frist / Slot_2h_frist_price
$ws.Range("AC2").Value = "8h-10h"
$ws.Range("AD2").Value = 86.43000000000001

# Slot_2h_second / Slot_2h_second_price
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 85.08

# Slot_min_price
$ws.Range("AG2").Value = "1h-23h"
